# Insert a new data row at row 194 (shifting existing rows 194:244 down to 195:245)
# and populate the new row with a new price record for "Primera" quality dated 2022-10-07.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 194; this pushes rows 194-244 to 195-245
# and automatically grows the sheet dimension to A1:T245.
$ws.Rows(194).Insert()

# Populate the newly inserted row 194 with the new record.
$ws.Cells.Item(194, 1).Value = 3
$ws.Cells.Item(194, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(194, 3).Value = "Coquimbo"
$ws.Cells.Item(194, 4).Value = 44841
$ws.Cells.Item(194, 5).Value = 5
$ws.Cells.Item(194, 6).Value = "Fruta"
$ws.Cells.Item(194, 7).Value = 100107
$ws.Cells.Item(194, 8).Value = "Otros"
$ws.Cells.Item(194, 9).Value = 100107002
$ws.Cells.Item(194, 10).Value = "Chirimoya"
$ws.Cells.Item(194, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(194, 12).Value = "Primera"
$ws.Cells.Item(194, 13).Value = 50
$ws.Cells.Item(194, 14).Value = 27000
$ws.Cells.Item(194, 15).Value = 27000
$ws.Cells.Item(194, 16).Value = 27000
$ws.Cells.Item(194, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(194, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(194, 19).Value = 2700
$ws.Cells.Item(194, 20).Value = 10
